# Riesgos workbook update:
#  - rename sheets back to default "Hoja1"/"Hoja2"
#  - fix "Riesgos Tecnicos" consequence (Media -> Alta)
#  - add a new "Riesgos Naturales" row to the risk list + resize the table
#  - add a stray space value further down the sheet
#  - extend / add the conditional-formatting rules that cover the new row
#  - restore the remembered selections on each sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Lista de Riesgos")
$ws2 = $wb.Worksheets.Item("Matriz de Riesgos")

# --- rename sheets -------------------------------------------------------
$ws1.Name = "Hoja1"
$ws2.Name = "Hoja2"

# --- fix existing data ----------------------------------------------------
$ws1.Range("C3").Value = "Alta"

# --- add the new risk row --------------------------------------------------
# copy formatting from an existing data row so the new cells pick up the
# same style used by the rest of the table
$ws1.Range("B3:D3").Copy()
$ws1.Range("B7:D7").PasteSpecial(-4122)
$ws1.Range("A7").Value = "Riesgos Naturales"
$ws1.Range("B7").Value = "Baja"
$ws1.Range("C7").Value = "Media"
$ws1.Range("D7").Value = 1

# stray cell further below the table
$ws1.Range("D12").Value = " "

# --- resize the structured table to include the new row --------------------
$lo = $ws1.ListObjects.Item("Tabla1")
$lo.Resize($ws1.Range("A1:D7"))

# --- conditional formatting --------------------------------------------------
# extend the range-wide rule that previously stopped at D6
$wideRules = $ws1.Range("D2:D6").FormatConditions
$wideRules.Item(1).ModifyAppliesToRange($ws1.Range("D2:D7"))

# add the rule for the new row, inserted as highest priority (mirrors how
# Excel slots a freshly authored rule ahead of the older ones)
$newRule = $ws1.Range("D7").FormatConditions.Add(2, $null, "AND(D7>=1*(D7>=1))")
$newRule.SetFirstPriority()

# --- restore remembered selections -----------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("C25").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B9").Select() | Out-Null
